$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 headers - reordered
$ws.Range("B1").Value = "Pint"
$ws.Range("C1").Value = "CABEAN"
$ws.Range("D1").Value = "PyBoolNet_Asp"
$ws.Range("E1").Value = "stable_motifs_new"
$ws.Range("F1").Value = "bioLQM"
$ws.Range("G1").Value = "boolsim"
$ws.Range("H1").Value = "sm_jgtz"
$ws.Range("I1").Value = "PyBoolNet"

# Row 2 values - new timing results
$ws.Range("B2").Value = 0.1798734390013124
$ws.Range("C2").Value = 0.07706998499998008
$ws.Range("D2").Value = 0.09504505499899096
$ws.Range("E2").Value = 0.6167890689994238
$ws.Range("F2").Value = 0.127908569000283
$ws.Range("G2").Value = 0.1977810940006748
$ws.Range("H2").Value = 31.96101380100117
$ws.Range("I2").Value = 0.3814105519995792
